# class/三年三班.xlsx - add two new students ("范总", "杨总") to Sheet1.
# Header row (row 1): A=name, B=pref2, C=pref1, D=pref3, E=sex, F=age, G=grade

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6: 范总
$ws.Range("A6").Value = "范总"
$ws.Range("B6").Value = "张三"
$ws.Range("C6").Value = "王柳"
$ws.Range("D6").Value = "王五"
$ws.Range("E6").Value = 0
$ws.Range("F6").Value = 23
$ws.Range("G6").Value = 88

# Row 7: 杨总
$ws.Range("A7").Value = "杨总"
$ws.Range("B7").Value = "王柳"
$ws.Range("C7").Value = "王五"
$ws.Range("D7").Value = "张三"
$ws.Range("E7").Value = 0
$ws.Range("F7").Value = 23
$ws.Range("G7").Value = 99

# Match the author's final cursor position recorded in the saved file.
$ws.Range("F12").Select()
